# App site's name changed and runners class fixed
# -> adds 20 new test-run rows (r=99..118) to the ScenarioStatus sheet,
#    each stamped with the new run date "27.12.22".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Creating a country", "FAILED", "chrome", "27.12.22"),
    @("Creating a country  with parameter data", "FAILED", "chrome", "27.12.22"),
    @("Login with valid username and password", "FAILED", "chrome", "27.12.22"),
    @("Login with valid username and password", "FAILED", "chrome", "27.12.22"),
    @("Login with valid username and password", "FAILED", "chrome", "27.12.22"),
    @("Login with valid username and password", "PASSED", "chrome", "27.12.22"),
    @("Creating a country", "PASSED", "chrome", "27.12.22"),
    @("Creating a country  with parameter data", "PASSED", "chrome", "27.12.22"),
    @("Creating a country", "PASSED", "chrome", "27.12.22"),
    @("Creating a country  with parameter data", "FAILED", "chrome", "27.12.22"),
    @("Creating a country", "PASSED", "firefox", "27.12.22"),
    @("Creating a country  with parameter data", "PASSED", "firefox", "27.12.22"),
    @("Creating a country", "PASSED", "chrome", "27.12.22"),
    @("Creating a country", "PASSED", "firefox", "27.12.22"),
    @("Creating a country  with parameter data", "FAILED", "chrome", "27.12.22"),
    @("Creating a country  with parameter data", "FAILED", "firefox", "27.12.22"),
    @("Creating a country", "PASSED", "chrome", "27.12.22"),
    @("Creating a country", "PASSED", "firefox", "27.12.22"),
    @("Creating a country  with parameter data", "PASSED", "chrome", "27.12.22"),
    @("Creating a country  with parameter data", "PASSED", "firefox", "27.12.22")
)

$startRow = 99
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
